$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.869288444519043
$ws.Range("B1").Value = 1.187809586524963
$ws.Range("C1").Value = 2.001517534255981
$ws.Range("D1").Value = 4.656040191650391
$ws.Range("E1").Value = 2.207576274871826
